$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pairs of rows whose data (columns B..AD, i.e. columns 2..30) must be swapped,
# keeping column A (row id) untouched.
$pairs = @(
    @(20, 21),
    @(22, 23),
    @(36, 37)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    for ($col = 2; $col -le 30; $col++) {
        $cell1 = $ws.Cells.Item($r1, $col)
        $cell2 = $ws.Cells.Item($r2, $col)

        $v1 = $cell1.Value2
        $v2 = $cell2.Value2

        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}
